$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.944.92"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.50%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.580.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "584.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.40"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.80%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -2.08%  "
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("E10").Value = "  -1.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.55%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.350"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.65%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.042.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "62.874.01"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("E16").Value = "  -2.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.585.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.10"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.32%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "341.84"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E21").Value = "  -3.52%  "
$ws.Range("E22").Value = "  +0.05%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.64"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.44%  "
$ws.Range("E24").Value = "  +6.08%  "
$ws.Range("E25").Value = "  -2.28%  "
$ws.Range("E26").Value = "  -3.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.24"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.36%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "460.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "0.0₃0800"
$ws.Range("E32").Value = "  -3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "176.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  +0.08%  "
$ws.Range("E36").Value = "  -2.02%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "18.85"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.52"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.02%  "
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("E40").Value = "  -3.18%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "158.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.51%  "
$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.95"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.36%  "
$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.635"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.26%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.35%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0538"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.50%  "
$ws.Range("E47").Value = "  -2.13%  "
$ws.Range("E48").Value = "  -1.45%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.68"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.50%  "
